$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 1.73
$ws.Range("G2").Value = 1.74
$ws.Range("H2").Value = 4.7
$ws.Range("Q2").Value = 1.46
$ws.Range("S2").Value = 2.1
$ws.Range("U2").Value = 2.74
$ws.Range("W2").Value = 2.34
$ws.Range("Y2").Value = 29
$ws.Range("AA2").Value = 95
$ws.Range("AB2").Value = 15.5
$ws.Range("AJ2").Value = 19.5
$ws.Range("AK2").Value = 15.5
$ws.Range("AM2").Value = 60
$ws.Range("AN2").Value = 6.2
$ws.Range("AO2").Value = 30

$ws.Range("F3").Value = 1.91
$ws.Range("G3").Value = 1.93
$ws.Range("H3").Value = 4.2
$ws.Range("S3").Value = 2.62
$ws.Range("W3").Value = 2.06
$ws.Range("AC3").Value = 10
$ws.Range("AK3").Value = 18.5

$ws.Range("Q4").Value = 1.57
$ws.Range("T4").Value = 1.66
$ws.Range("U4").Value = 2.22

$ws.Range("M5").Value = 1.06
$ws.Range("R5").Value = 1.09

$ws.Range("M6").Value = 1.05
$ws.Range("O6").Value = 1.06

$ws.Range("F7").Value = 14.5
$ws.Range("I7").Value = 1.27
$ws.Range("J7").Value = 7.2

$ws.Range("M8").Value = 1.05
$ws.Range("O8").Value = 1.06
$ws.Range("Q8").Value = 1.37

$ws.Range("P9").Value = 1.84
$ws.Range("U9").Value = 2.06
$ws.Range("AH9").Value = 18.5

$ws.Range("F10").Value = 1.57
$ws.Range("G10").Value = 1.58
$ws.Range("Q10").Value = 1.65
$ws.Range("V10").Value = 1.17
$ws.Range("W10").Value = 2.72
$ws.Range("X10").Value = 22

$ws.Range("K11").Value = 4.6
$ws.Range("L11").Value = 1.2
$ws.Range("O11").Value = 1.18
$ws.Range("P11").Value = 2.52
$ws.Range("Q11").Value = 1.53
$ws.Range("R11").Value = 1.61
$ws.Range("S11").Value = 2.3
$ws.Range("T11").Value = 1.54

$ws.Range("H12").Value = 3.3
$ws.Range("O12").Value = 1.32
$ws.Range("U12").Value = 2.2
$ws.Range("AJ12").Value = 32

$ws.Range("F13").Value = 8.8
$ws.Range("G13").Value = 9
$ws.Range("U13").Value = 2.02
$ws.Range("AF13").Value = 80
$ws.Range("AH13").Value = 24
$ws.Range("AO13").Value = 5.8
